$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.811.99'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.547.55'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '310.26'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -2.13%  '
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '98.51'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +1.34%  '
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.38%  '
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.531'
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -0.91%  '
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '35.73'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -1.43%  '
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0807'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.40'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -2.94%  '
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.941.25'
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +0.40%  '
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.89'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +4.66%  '
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.632.36'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +3.43%  '
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.840'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '42.833.25'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.32%  '
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.76'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -1.80%  '
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.43'
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -2.84%  '
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0957'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.25%  '
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '69.34'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '248.20'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -2.26%  '
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.92'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -1.34%  '
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +0.37%  '
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '26.69'
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.70%  '
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -1.62%  '
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '40.32'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.69%  '
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '10.14'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -2.97%  '
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '158.01'
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.20%  '
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.76'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -2.88%  '
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0801'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.29'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -2.58%  '
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.09'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -3.95%  '
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.63'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -3.53%  '
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +4.34%  '
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '18.38'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.03%  '
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -1.59%  '
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +2.22%  '
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '4.10'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +6.66%  '
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0300'
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -1.38%  '
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.000.51'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -1.15%  '
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -2.44%  '
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.06'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -0.53%  '
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.794.27'
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +0.34%  '
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '81.43'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -3.52%  '
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +0.31%  '
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '73.67'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -4.18%  '
$c.ClearFormats()
